$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.498833115296896
$ws.Range("B2").Value = 0.501166884703103

$ws.Range("A3").Select()
